# New .ttl from Google sheet has been generated.
# A new "dct:creator" entry (Eva, with her ORCID) is inserted above the
# existing "dct:creator" (Hannah) row, pushing all subsequent metadata /
# term rows down by one row.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row above row 12 (dct:creator / Hannah), shifting
# everything below (rows 12-89) down by one (to rows 13-90).
$ws.Rows("12:12").Insert()

# Populate the newly inserted row with the new creator's details.
$ws.Range("A12").Value = "dct:creator"
$ws.Range("B12").Value = "https://orcid.org/0000-0003-4093-2147"
$ws.Range("C12").Value = "Eva"
